# Add a new paragraph "Added new line" right after the paragraph that
# contains "Using powershell." (the last paragraph in the body), matching
# its run/paragraph formatting (en-US language).

$d = $word.ActiveDocument

# Locate the paragraph to insert after by searching for its text so the
# script is resilient to exact paragraph indices.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Using powershell.*") {
        $target = $candidate
    }
}

if ($target -eq $null) {
    # Fallback: just use the last paragraph in the document.
    $target = $d.Paragraphs.Item($d.Paragraphs.Count)
}

$r = $target.Range
$r.Collapse(0)
$r.InsertParagraphAfter()

# The freshly inserted (now-last) paragraph is empty; fill it with the text.
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.InsertAfter("Added new line")

Write-Output "Inserted new paragraph with 'Added new line'."
